$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the long shell-script cell (C2): reorder DEBIAN_FRONTEND export + add
# a dpkg conf-file prompt override, and drop the do-release-upgrade override flags.
$ws.Range("C2").Value = "echo 'Dpkg::Options { `"--force-confdef`"; `"--force-confold`"; }' | sudo tee /etc/apt/apt.conf.d/99force-confold`nexport DEBIAN_FRONTEND=noninteractive`nsudo apt-get update`nsudo apt-get upgrade -y`nsudo apt-get dist-upgrade -y`nsudo apt-get autoremove -y`nsudo apt-get install -y update-manager-core`nsudo sed -i 's/^Prompt=lts/Prompt=normal/' /etc/update-manager/release-upgrades`nsudo -E do-release-upgrade -f DistUpgradeViewNonInteractive -m server`ncat /etc/os-release"

# Row 2 was manually shrunk from the Excel max (409.5) down to 174 points.
$ws.Rows.Item(2).RowHeight = 174

# Selection moved from D2 to B2.
$ws.Range("B2").Select()
